# Add a new "Swiss" test-data worksheet, cloned from the existing "Czech"
# sheet (same layout/formatting), then populate it with the Switzerland
# market values. Also update the previously-selected "Czech" sheet so the
# new "Swiss" sheet becomes the active / selected tab, matching how Excel
# leaves selection state after adding+activating a new sheet.

$wb = $excel.ActiveWorkbook

$czech = $wb.Worksheets.Item("Czech")

# Clear the Czech sheet's previous selection state (it was the tab that
# used to be active/selected) before handing focus to the new sheet.
$czech.Cells.Select() | Out-Null

# Duplicate the Czech sheet (keeps columns widths, styles, merged cells,
# page setup, etc.) and place the copy right after it.
$czech.Copy($null, $czech)

# The copy becomes the last sheet and is activated automatically.
$swiss = $wb.Worksheets.Item($wb.Worksheets.Count)
$swiss.Name = "Swiss"

# Populate the Switzerland-specific values.
$swiss.Range("A1").Value = "S"
$swiss.Range("B2").Value = "Switzerland Market"
$swiss.Range("B4").Value = "NGC-3476/T2352"

# Match the selection shown on the new sheet.
$swiss.Range("B2:B4").Select() | Out-Null
